$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '74.830.96'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '2.812.16'
$ws.Range("E3").Value = '  +6.83%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'187.75"
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").Value = "'591.36"
$ws.Range("E6").Value = '  +1.37%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +2.53%  '
$ws.Range("E9").Value = '  -4.96%  '
$ws.Range("D10").Value = '2.810.80'
$ws.Range("E10").Value = '  +6.76%  '
$ws.Range("E11").Value = '  +5.23%  '
$ws.Range("E12").Value = '  -2.01%  '
$ws.Range("E13").Value = '  +2.98%  '
$ws.Range("D14").Value = '3.330.97'
$ws.Range("E14").Value = '  +6.89%  '
$ws.Range("D15").Value = '74.805.66'
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("E16").Value = '  -1.51%  '
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").Value = '2.813.89'
$ws.Range("E18").Value = '  +7.16%  '
$ws.Range("D19").Value = "'8.88"
$ws.Range("E19").Value = '  -3.85%  '
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("D21").Value = "'376.29"
$ws.Range("E21").Value = '  +2.83%  '
$ws.Range("D22").Value = "'2.27"
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("D23").Value = "'4.08"
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = "'70.72"
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("D26").Value = '2.952.10'
$ws.Range("E26").Value = '  +7.06%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").Value = "'9.67"
$ws.Range("E28").Value = '  +3.78%  '
$ws.Range("E29").Value = '  +8.63%  '
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").Value = "'510.33"
$ws.Range("E32").Value = '  -2.74%  '
$ws.Range("D33").Value = "'7.59"
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = "'1.79"
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").Value = "'164.31"
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("E37").Value = '  +3.77%  '
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("D39").Value = "'19.36"
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("D40").Value = "'181.74"
$ws.Range("E40").Value = '  +11.41%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +4.56%  '
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = "'40.05"
$ws.Range("E45").Value = '  +2.73%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").Value = "'1.20"
$ws.Range("E46").Value = '  +2.12%  '
$ws.Range("D47").Value = "'0.0862"
$ws.Range("E47").Value = '  +2.05%  '
$ws.Range("E48").Value = '  -3.38%  '
$ws.Range("D49").Value = "'0.572"
$ws.Range("E49").Value = '  +9.18%  '
$ws.Range("E50").Value = '  +2.96%  '
$ws.Range("D51").Value = "'0.633"
$ws.Range("E51").Value = '  +7.65%  '
